$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Cells.Item(2, 11).Value = "2025-10-27T16:21:29.502651+00:00"
$ws.Cells.Item(3, 11).Value = "2025-10-27T16:21:31.587520+00:00"
$ws.Cells.Item(4, 11).Value = "2025-10-27T16:21:31.587540+00:00"
$ws.Cells.Item(5, 11).Value = "2025-10-27T16:21:31.587550+00:00"
$ws.Cells.Item(6, 11).Value = "2025-10-27T16:21:31.587558+00:00"
$ws.Cells.Item(7, 11).Value = "2025-10-27T16:21:31.587566+00:00"
$ws.Cells.Item(8, 11).Value = "2025-10-27T16:21:33.679237+00:00"
$ws.Cells.Item(9, 11).Value = "2025-10-27T16:21:33.679255+00:00"
$ws.Cells.Item(10, 11).Value = "2025-10-27T16:21:35.685236+00:00"
$ws.Cells.Item(11, 11).Value = "2025-10-27T16:21:38.537677+00:00"
$ws.Cells.Item(12, 11).Value = "2025-10-27T16:21:38.537693+00:00"
$ws.Cells.Item(13, 11).Value = "2025-10-27T16:21:40.700179+00:00"
$ws.Cells.Item(14, 11).Value = "2025-10-27T16:21:40.700196+00:00"
$ws.Cells.Item(15, 11).Value = "2025-10-27T16:21:40.700204+00:00"
$ws.Cells.Item(16, 11).Value = "2025-10-27T16:21:40.700212+00:00"
$ws.Cells.Item(17, 11).Value = "2025-10-27T16:21:47.143774+00:00"
$ws.Cells.Item(18, 11).Value = "2025-10-27T16:21:49.185950+00:00"
$ws.Cells.Item(19, 11).Value = "2025-10-27T16:21:51.176975+00:00"
$ws.Cells.Item(20, 11).Value = "2025-10-27T16:21:53.231308+00:00"
$ws.Cells.Item(21, 11).Value = "2025-10-27T16:21:53.231327+00:00"
$ws.Cells.Item(22, 11).Value = "2025-10-27T16:21:55.667200+00:00"
$ws.Cells.Item(23, 11).Value = "2025-10-27T16:21:55.667226+00:00"
$ws.Cells.Item(24, 11).Value = "2025-10-27T16:21:55.667234+00:00"
$ws.Cells.Item(25, 11).Value = "2025-10-27T16:21:55.667241+00:00"
$ws.Cells.Item(26, 11).Value = "2025-10-27T16:21:55.667249+00:00"
$ws.Cells.Item(27, 11).Value = "2025-10-27T16:21:58.082418+00:00"
$ws.Cells.Item(28, 11).Value = "2025-10-27T16:22:02.551092+00:00"
$ws.Cells.Item(29, 11).Value = "2025-10-27T16:22:02.551110+00:00"
$ws.Cells.Item(30, 11).Value = "2025-10-27T16:22:02.551117+00:00"
$ws.Cells.Item(31, 11).Value = "2025-10-27T16:22:02.551125+00:00"
$ws.Cells.Item(32, 11).Value = "2025-10-27T16:22:04.585314+00:00"
$ws.Cells.Item(33, 11).Value = "2025-10-27T16:22:04.585331+00:00"
$ws.Cells.Item(34, 11).Value = "2025-10-27T16:22:06.506160+00:00"
$ws.Cells.Item(35, 11).Value = "2025-10-27T16:22:06.506177+00:00"
$ws.Cells.Item(36, 11).Value = "2025-10-27T16:22:06.506184+00:00"
$ws.Cells.Item(37, 11).Value = "2025-10-27T16:22:06.506191+00:00"
$ws.Cells.Item(38, 11).Value = "2025-10-27T16:22:06.506198+00:00"
$ws.Cells.Item(39, 11).Value = "2025-10-27T16:22:06.506206+00:00"
$ws.Cells.Item(40, 11).Value = "2025-10-27T16:22:06.506213+00:00"
$ws.Cells.Item(41, 11).Value = "2025-10-27T16:22:06.506220+00:00"
$ws.Cells.Item(42, 11).Value = "2025-10-27T16:22:06.506226+00:00"
$ws.Cells.Item(43, 11).Value = "2025-10-27T16:22:08.584087+00:00"
$ws.Cells.Item(44, 11).Value = "2025-10-27T16:22:08.584104+00:00"
$ws.Cells.Item(45, 11).Value = "2025-10-27T16:22:13.103979+00:00"
$ws.Cells.Item(46, 11).Value = "2025-10-27T16:22:15.101542+00:00"
$ws.Cells.Item(47, 11).Value = "2025-10-27T16:22:15.101559+00:00"
$ws.Cells.Item(48, 11).Value = "2025-10-27T16:22:15.101567+00:00"
$ws.Cells.Item(49, 11).Value = "2025-10-27T16:22:15.101574+00:00"

$wb.Save()
